$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Septiembre de 2020 a las 22:24"

# Update country rows whose stats / ranking changed (re-sorted by Casos totales descending)
# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 7124958
$ws.Range("C4").Value = 27021
$ws.Range("D4").Value = 4378363
$ws.Range("E4").Value = 2540344
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 781
$ws.Range("H4").Value = 206251
# Row 5: India
$ws.Range("A5").Value = "India"
$ws.Range("B5").Value = 5730184
$ws.Range("C5").Value = 89688
$ws.Range("D5").Value = 4671850
$ws.Range("E5").Value = 967161
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1152
$ws.Range("H5").Value = 91173
# Row 12: Sudafrica
$ws.Range("A12").Value = "Sudafrica"
$ws.Range("B12").Value = 665188
$ws.Range("C12").Value = 1906
$ws.Range("D12").Value = 594229
$ws.Range("E12").Value = 54753
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 88
$ws.Range("H12").Value = 16206
# Row 25: Alemania
$ws.Range("A25").Value = "Alemania"
$ws.Range("B25").Value = 279160
$ws.Range("C25").Value = 1984
$ws.Range("D25").Value = 247900
$ws.Range("E25").Value = 21752
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 9508
# Row 29: Canada
$ws.Range("A29").Value = "Canada"
$ws.Range("B29").Value = 147522
$ws.Range("C29").Value = 859
$ws.Range("D29").Value = 127422
$ws.Range("E29").Value = 10858
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = 9242
# Row 54: Costa Rica
$ws.Range("A54").Value = "Costa Rica"
$ws.Range("B54").Value = 68059
$ws.Range("C54").Value = 1370
$ws.Range("D54").Value = 26136
$ws.Range("E54").Value = 41142
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 21
$ws.Range("H54").Value = 781
# Row 55: Nepal
$ws.Range("A55").Value = "Nepal"
$ws.Range("B55").Value = 67804
$ws.Range("C55").Value = 1172
$ws.Range("D55").Value = 49954
$ws.Range("E55").Value = 17414
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 7
$ws.Range("H55").Value = 436
# Row 84: Costa de Marfil
$ws.Range("A84").Value = "Costa de Marfil"
$ws.Range("B84").Value = 19430
$ws.Range("C84").Value = 87
$ws.Range("D84").Value = 18875
$ws.Range("E84").Value = 435
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 120
# Row 109: Mauritania
$ws.Range("A109").Value = "Mauritania"
$ws.Range("B109").Value = 7425
$ws.Range("C109").Value = 22
$ws.Range("D109").Value = 7028
$ws.Range("E109").Value = 236
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 161
# Row 115: Malaui
$ws.Range("A115").Value = "Malaui"
$ws.Range("B115").Value = 5746
$ws.Range("C115").Value = 7
$ws.Range("D115").Value = 4140
$ws.Range("E115").Value = 1427
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 179
# Row 126: Ruanda
$ws.Range("A126").Value = "Ruanda"
$ws.Range("B126").Value = 4779
$ws.Range("C126").Value = 41
$ws.Range("D126").Value = 2995
$ws.Range("E126").Value = 1757
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 27
# Row 127: Surinam
$ws.Range("A127").Value = "Surinam"
$ws.Range("B127").Value = 4759
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 4541
$ws.Range("E127").Value = 118
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 100
# Row 143: Mali
$ws.Range("A143").Value = "Mali"
$ws.Range("B143").Value = 3034
$ws.Range("C143").Value = 4
$ws.Range("D143").Value = 2382
$ws.Range("E143").Value = 522
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 130
# Row 144: Estonia
$ws.Range("A144").Value = "Estonia"
$ws.Range("B144").Value = 3033
$ws.Range("C144").Value = 57
$ws.Range("D144").Value = 2387
$ws.Range("E144").Value = 582
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 64
# Row 152: Sierra Leona
$ws.Range("A152").Value = "Sierra Leona"
$ws.Range("B152").Value = 2183
$ws.Range("C152").Value = 9
$ws.Range("D152").Value = 1665
$ws.Range("E152").Value = 446
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 72
# Row 157: Principado de Andorra
$ws.Range("A157").Value = "Principado de Andorra"
$ws.Range("B157").Value = 1753
$ws.Range("C157").Value = 72
$ws.Range("D157").Value = 1203
$ws.Range("E157").Value = 497
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 53
# Row 158: Togo
$ws.Range("A158").Value = "Togo"
$ws.Range("B158").Value = 1701
$ws.Range("C158").Value = 18
$ws.Range("D158").Value = 1297
$ws.Range("E158").Value = 363
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 41
# Row 167: Republica del Chad
$ws.Range("A167").Value = "Republica del Chad"
$ws.Range("B167").Value = 1164
$ws.Range("C167").Value = 9
$ws.Range("D167").Value = 997
$ws.Range("E167").Value = 85
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 1
$ws.Range("H167").Value = 82
# Row 214: Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
# Row 215: Islas Malvinas
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
